$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (rich-text shared strings: Volume/Number and report week dates) ---
$ws.Range("A8").Value = "Volume 31   Number  5"
$ws.Range("C9").Value = "Report Covering the Week  1/29/2024  Through  2/4/2024"

# --- Cells that become "N/A"-style text placeholders: copy value+style from a cell that
#     already holds the same shared-string/style combination ---
$ws.Range("C14").Copy($ws.Range("C18"))
$ws.Range("D14").Copy($ws.Range("D17"))
$ws.Range("F14").Copy($ws.Range("D27"))
$ws.Range("L14").Copy($ws.Range("E17"))
$ws.Range("M27").Copy($ws.Range("E27"))

# --- Cells changing FROM text-placeholder style TO a numeric style: set the number format
#     from a same-row cell that already has the destination style, then set the value ---
$ws.Range("M14").NumberFormat = $ws.Range("N14").NumberFormat
$ws.Range("M14").Value = -100
$ws.Range("C15").NumberFormat = $ws.Range("G15").NumberFormat
$ws.Range("C15").Value = 1
$ws.Range("L15").NumberFormat = $ws.Range("H15").NumberFormat
$ws.Range("L15").Value = 200
$ws.Range("D23").NumberFormat = $ws.Range("C23").NumberFormat
$ws.Range("D23").Value = 4
$ws.Range("E23").NumberFormat = $ws.Range("H23").NumberFormat
$ws.Range("E23").Value = -50
$ws.Range("C26").NumberFormat = $ws.Range("G26").NumberFormat
$ws.Range("C26").Value = 1
$ws.Range("L26").NumberFormat = $ws.Range("H26").NumberFormat
$ws.Range("L26").Value = 200

# --- Plain numeric value updates (style unchanged) ---
$ws.Range("F15").Value = 3
$ws.Range("H15").Value = 200
$ws.Range("I15").Value = 3
$ws.Range("K15").Value = 200
$ws.Range("D16").Value = 4
$ws.Range("E16").Value = -25
$ws.Range("G16").Value = 15
$ws.Range("H16").Value = -33.333333333333
$ws.Range("I16").Value = 13
$ws.Range("J16").Value = 18
$ws.Range("K16").Value = -27.777777777777
$ws.Range("L16").Value = -27.777777777777
$ws.Range("M16").Value = -31.578947368421
$ws.Range("N16").Value = -80.882352941176
$ws.Range("F17").Value = 7
$ws.Range("G17").Value = 7
$ws.Range("H17").Value = 0
$ws.Range("I17").Value = 13
$ws.Range("K17").Value = 30
$ws.Range("L17").Value = -7.142857142857
$ws.Range("M17").Value = 0
$ws.Range("N17").Value = -65.789473684210
$ws.Range("F18").Value = 11
$ws.Range("G18").Value = 5
$ws.Range("H18").Value = 120
$ws.Range("L18").Value = -50
$ws.Range("M18").Value = -23.529411764705
$ws.Range("N18").Value = -89.6
$ws.Range("C19").Value = 9
$ws.Range("D19").Value = 8
$ws.Range("E19").Value = 12.5
$ws.Range("F19").Value = 32
$ws.Range("H19").Value = -33.333333333333
$ws.Range("I19").Value = 38
$ws.Range("J19").Value = 56
$ws.Range("K19").Value = -32.142857142857
$ws.Range("L19").Value = -22.448979591836
$ws.Range("M19").Value = -13.636363636363
$ws.Range("N19").Value = -61.616161616161
$ws.Range("C20").Value = 1
$ws.Range("D20").Value = 3
$ws.Range("E20").Value = -66.666666666666
$ws.Range("F20").Value = 4
$ws.Range("G20").Value = 12
$ws.Range("H20").Value = -66.666666666666
$ws.Range("I20").Value = 4
$ws.Range("J20").Value = 14
$ws.Range("K20").Value = -71.428571428571
$ws.Range("L20").Value = 0
$ws.Range("M20").Value = 33.333333333333
$ws.Range("N20").Value = -96.039603960396
$ws.Range("C21").Value = 17
$ws.Range("D21").Value = 15
$ws.Range("E21").Value = 13.333333333333
$ws.Range("F21").Value = 67
$ws.Range("G21").Value = 89
$ws.Range("H21").Value = -24.719101123595
$ws.Range("I21").Value = 84
$ws.Range("J21").Value = 110
$ws.Range("K21").Value = -23.636363636363
$ws.Range("L21").Value = -25
$ws.Range("M21").Value = -13.402061855670
$ws.Range("N21").Value = -80.821917808219
$ws.Range("C22").Value = 2
$ws.Range("E22").Value = 100
$ws.Range("F22").Value = 4
$ws.Range("G22").Value = 3
$ws.Range("H22").Value = 33.333333333333
$ws.Range("I22").Value = 5
$ws.Range("J22").Value = 3
$ws.Range("K22").Value = 66.666666666666
$ws.Range("L22").Value = 25
$ws.Range("M22").Value = 66.666666666666
$ws.Range("F23").Value = 4
$ws.Range("G23").Value = 7
$ws.Range("H23").Value = -42.857142857142
$ws.Range("I23").Value = 8
$ws.Range("J23").Value = 8
$ws.Range("K23").Value = 0
$ws.Range("L23").Value = -20
$ws.Range("M23").Value = -11.111111111111
$ws.Range("C24").Value = 22
$ws.Range("D24").Value = 44
$ws.Range("E24").Value = -50
$ws.Range("F24").Value = 100
$ws.Range("G24").Value = 196
$ws.Range("H24").Value = -48.979591836734
$ws.Range("I24").Value = 132
$ws.Range("J24").Value = 236
$ws.Range("K24").Value = -44.067796610169
$ws.Range("L24").Value = -17.5
$ws.Range("M24").Value = 25.714285714285
$ws.Range("C25").Value = 6
$ws.Range("D25").Value = 4
$ws.Range("E25").Value = 50
$ws.Range("F25").Value = 22
$ws.Range("H25").Value = -8.333333333333
$ws.Range("I25").Value = 31
$ws.Range("J25").Value = 28
$ws.Range("K25").Value = 10.714285714285
$ws.Range("L25").Value = 0
$ws.Range("M25").Value = 10.714285714285
$ws.Range("F26").Value = 3
$ws.Range("H26").Value = 200
$ws.Range("I26").Value = 3
$ws.Range("K26").Value = 200
$ws.Range("F27").Value = 1
$ws.Range("G27").Value = 2
$ws.Range("H27").Value = -50
$ws.Range("L27").Value = -50
$ws.Range("F30").Value = 2
$ws.Range("I30").Value = 3
